# Rename the worksheet: "Property1" -> "DataNode"
# (reflects the commit's unification of the DataNode/DataTable/Entity concepts)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move the active selection in the frozen (bottom-left) pane from D9 to F25
$win = $excel.ActiveWindow
$win.Panes.Item(2).Activate()
$ws.Range("F25").Select()
